$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "eth"
$ws.Range("B35").Value = "Correlated Concept"
$ws.Range("C35").Value = "ETC"
$ws.Range("D35").Value = "test111"

$ws.Range("A36").Value = "eth"
$ws.Range("B36").Value = "Correlated Concept"
$ws.Range("C36").Value = "DeFi"
$ws.Range("D36").Value = "test111"
